$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 797. This shifts the existing rows 797..894
# down to 798..895, matching the dimension change from A1:R894 to A1:R895.
$ws.Rows.Item(797).Insert()

# Populate the newly inserted row 797 with its values.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are identical to the row that used to be
# at 797 (now shifted to 798); only D,J,K,L,M,P differ for the new record.
$ws.Cells.Item(797, 1).Value = 3
$ws.Cells.Item(797, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(797, 3).Value = "Coquimbo"
$ws.Cells.Item(797, 4).Value = 45212
$ws.Cells.Item(797, 5).Value = 5
$ws.Cells.Item(797, 6).Value = 100112003
$ws.Cells.Item(797, 7).Value = "Ajo"
$ws.Cells.Item(797, 8).Value = "Chino"
$ws.Cells.Item(797, 9).Value = "Primera"
$ws.Cells.Item(797, 10).Value = 40
$ws.Cells.Item(797, 11).Value = 19000
$ws.Cells.Item(797, 12).Value = 19000
$ws.Cells.Item(797, 13).Value = 19000
$ws.Cells.Item(797, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(797, 15).Value = "China"
$ws.Cells.Item(797, 16).Value = 1900
$ws.Cells.Item(797, 17).Value = 10
$ws.Cells.Item(797, 18).Value = "Hortaliza"
